# Add 20% growth rate constraints for PV and CSP on the PWR sheet.
#
# The PWR sheet has two related blocks:
#   - rows 6-11 : "UC" constraint rows, each built from a matching data row
#                 further down the sheet (columns A-D) via formulas.
#   - rows 14-19: the underlying data rows (Technology/Unit/Growth/Decline).
#
# We insert two new rows in each block (mirroring the existing CCS/Wave/
# Tidal/WON/WOF pattern) and populate them for the new PV and CSP
# technologies.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PWR")

# --- Insert two new rows in the upper "UC" block (after row 11, the WOF
#     row), pushing the blank row 13 + data block down by two rows. ---
$ws.Rows("12:13").Insert()

# --- Insert two new rows in the lower data block (after row 21, the WOF
#     data row) for the new PV and CSP entries. ---
$ws.Rows("22:23").Insert()

# --- Populate the new data rows 22 (PV) and 23 (CSP), mirroring rows
#     20/21 (WON/WOF). ---
$ws.Cells.Item(22, 1).Value = "PWR"
$ws.Cells.Item(22, 2).Value = "PV"
$ws.Cells.Item(22, 3).Value = 0.2
$ws.Cells.Item(22, 4).Value = 0.5

$ws.Cells.Item(23, 1).Value = "PWR"
$ws.Cells.Item(23, 2).Value = "CSP"
$ws.Cells.Item(23, 3).Value = 0.2
$ws.Cells.Item(23, 4).Value = 0.5

# --- Populate the new "UC" rows 12 (PV) and 13 (CSP), mirroring rows
#     10/11 (WON/WOF) but pointing at the new data rows 22/23. ---
$ws.Range("B12").Formula = '=_xlfn.TEXTJOIN("_",TRUE,"UC",A22,"MaxGrowth",B22)'
$ws.Range("C12").Value = "CAP, GROWTH"
$ws.Range("G12").Value = "P*SOL*PV*"
$ws.Cells.Item(12, 8).Value = 2018
$ws.Range("I12").Value = "LO"
$ws.Range("J12").Formula = "=1+`$C22"
$ws.Cells.Item(12, 11).Value = 1
$ws.Range("L12").Formula = "=-D22"
$ws.Cells.Item(12, 13).Value = 5
$ws.Range("N12").Formula = '=_xlfn.TEXTJOIN(" ",TRUE,A22, "maximum growth rate of",B22)'

$ws.Range("B13").Formula = '=_xlfn.TEXTJOIN("_",TRUE,"UC",A23,"MaxGrowth",B23)'
$ws.Range("C13").Value = "CAP, GROWTH"
$ws.Range("G13").Value = "P*SOL*CSP*"
$ws.Cells.Item(13, 8).Value = 2018
$ws.Range("I13").Value = "LO"
$ws.Range("J13").Formula = "=1+`$C23"
$ws.Cells.Item(13, 11).Value = 1
$ws.Range("L13").Formula = "=-D23"
$ws.Cells.Item(13, 13).Value = 5
$ws.Range("N13").Formula = '=_xlfn.TEXTJOIN(" ",TRUE,A23, "maximum growth rate of",B23)'

# Match the author's final selection (cell I13, the "LO" cell of the new
# CSP row) recorded in the sheet view.
$ws.Activate()
$ws.Range("I13").Select()
